$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: these "PickupID/POD No" values look numeric but must be
# stored as text (matches the source data, which comes from shared strings,
# not numbers). Briefly force a Text number format so Excel keeps the
# leading value as a string, then clear the format again so the cell keeps
# its original (default) style - only the underlying value/type changes.
function Set-TextValue($rangeAddress, $value) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "C11" "10266615"
Set-TextValue "C12" "10264902"
Set-TextValue "C13" "10264903"
Set-TextValue "C14" "10264909"
Set-TextValue "C24" "136895718"

$ws.Range("F26").Value = 'Cannot invoke "org.openqa.selenium.WebElement.isDisplayed()" because "element" is null'
